$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new mapped field names to column B (rows 9 and 10)
$ws.Range("B9").Value = "companies_id"
$ws.Range("B10").Value = "extra_hour_distributions_id"

# Widen column B to fit the new longer text
$ws.Columns.Item(2).ColumnWidth = 26.1666666666667

# Move the active selection to B11, matching the saved view state
$ws.Range("B11").Select()
